$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 356.72223
$ws.Range("I5").Value = 330.46155
$ws.Range("K5").Value = 330.46155
$ws.Range("M5").Value = -215.46155

$ws.Range("H51").Value = 2863.5557
$ws.Range("I51").Value = 885
$ws.Range("J51").Value = 3428.8572
$ws.Range("K51").Value = 885
$ws.Range("L51").Value = 3428.8572
$ws.Range("M51").Value = -401
$ws.Range("N51").Value = -4396.8572

$ws.Range("H125").Value = 2090
$ws.Range("I125").Value = 3316.6667
$ws.Range("K125").Value = 29850.0003
$ws.Range("M125").Value = -27390.0003

$ws.Range("H137").Value = 1583.2903
$ws.Range("I137").Value = 1316.1305
$ws.Range("J137").Value = 2351.375
$ws.Range("K137").Value = 3948.3915
$ws.Range("L137").Value = 7054.125
$ws.Range("M137").Value = -1398.3915
$ws.Range("N137").Value = -12154.125

$ws.Range("H138").Value = 532159.25
$ws.Range("I138").Value = 1982.5714
$ws.Range("J138").Value = 581642.4
$ws.Range("K138").Value = 5947.7142
$ws.Range("L138").Value = 1744927.2
$ws.Range("M138").Value = -807.7142000000003
$ws.Range("N138").Value = -1755207.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3141.83
$ws.Range("I32").Value = 3409.1777
$ws.Range("J32").Value = 1638
$ws.Range("K32").Value = 3409.1777
$ws.Range("L32").Value = 1638
$ws.Range("M32").Value = -3122.1777
$ws.Range("N32").Value = -2212

$ws.Range("H74").Value = 1170.1305
$ws.Range("I74").Value = 956.6111
$ws.Range("K74").Value = 956.6111
$ws.Range("M74").Value = -82.61109999999996

$ws.Range("H77").Value = 1170.1305
$ws.Range("I77").Value = 956.6111
$ws.Range("K77").Value = 4783.055499999999
$ws.Range("M77").Value = -415.0554999999995

$ws.Range("H135").Value = 17203.5
$ws.Range("J135").Value = 17203.5
$ws.Range("L135").Value = 17203.5
$ws.Range("N135").Value = -27343.5

$ws.Range("H139").Value = 48857.5
$ws.Range("J139").Value = 48857.5
$ws.Range("L139").Value = 48857.5
$ws.Range("N139").Value = -59137.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1553.9445
$ws.Range("J20").Value = 1788.5
$ws.Range("L20").Value = 1788.5
$ws.Range("N20").Value = -2282.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1426.375
$ws.Range("I31").Value = 1155.2222
$ws.Range("J31").Value = 1648.2273
$ws.Range("K31").Value = 1155.2222
$ws.Range("L31").Value = 1648.2273
$ws.Range("M31").Value = -860.2221999999999
$ws.Range("N31").Value = -2238.2273

$ws.Range("H34").Value = 1426.375
$ws.Range("I34").Value = 1155.2222
$ws.Range("J34").Value = 1648.2273
$ws.Range("K34").Value = 1155.2222
$ws.Range("L34").Value = 1648.2273
$ws.Range("M34").Value = -953.2221999999999
$ws.Range("N34").Value = -2052.2273

$ws.Range("H58").Value = 675
$ws.Range("I58").Value = 666.82355
$ws.Range("K58").Value = 666.82355
$ws.Range("M58").Value = -463.82355

$ws.Range("H94").Value = 858
$ws.Range("I94").Value = 2099.5
$ws.Range("J94").Value = 609.7
$ws.Range("K94").Value = 2099.5
$ws.Range("L94").Value = 609.7
$ws.Range("M94").Value = -1648.5
$ws.Range("N94").Value = -1511.7

$ws.Range("H134").Value = 7093343
$ws.Range("I134").Value = 9260328
$ws.Range("K134").Value = 27780984
$ws.Range("M134").Value = -27778449

$ws.Range("H136").Value = 675
$ws.Range("I136").Value = 666.82355
$ws.Range("K136").Value = 2000.47065
$ws.Range("M136").Value = 549.5293500000002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 474386.72
$ws.Range("I4").Value = 90066.39999999999
$ws.Range("J4").Value = 666546.9
$ws.Range("K4").Value = 270199.2
$ws.Range("L4").Value = 1999640.7
$ws.Range("M4").Value = -270087.2
$ws.Range("N4").Value = -1999864.7

$ws.Range("H23").Value = 203.63637
$ws.Range("I23").Value = 118
$ws.Range("J23").Value = 275
$ws.Range("K23").Value = 354
$ws.Range("L23").Value = 825
$ws.Range("M23").Value = -119
$ws.Range("N23").Value = -1295

$ws.Range("H34").Value = 3126710.2
$ws.Range("J34").Value = 3705617.8
$ws.Range("L34").Value = 11116853.4
$ws.Range("N34").Value = -11117021.4

$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()

$ws.Range("H55").Value = 3249.8333
$ws.Range("J55").Value = 3249.8333
$ws.Range("L55").Value = 9749.499899999999
$ws.Range("N55").Value = -10103.4999

$ws.Range("H68").Value = 1968.5193
$ws.Range("J68").Value = 1995.3529
$ws.Range("L68").Value = 5986.0587
$ws.Range("N68").Value = -7608.0587

$ws.Range("H71").Value = 1968.5193
$ws.Range("J71").Value = 1995.3529
$ws.Range("L71").Value = 17958.1761
$ws.Range("N71").Value = -26070.1761

$ws.Range("H107").Value = 7057.8125
$ws.Range("I107").Value = 675.5
$ws.Range("J107").Value = 10887.2
$ws.Range("K107").Value = 2026.5
$ws.Range("L107").Value = 32661.6
$ws.Range("M107").Value = -106.5
$ws.Range("N107").Value = -36501.60000000001

$ws.Range("H132").Value = 1066.1666
$ws.Range("I132").Value = 649.5
$ws.Range("J132").Value = 1899.5
$ws.Range("K132").Value = 5845.5
$ws.Range("L132").Value = 17095.5
$ws.Range("M132").Value = -3315.5
$ws.Range("N132").Value = -22155.5

$ws.Range("H140").Value = 33758.113
$ws.Range("I140").Value = 40129.242
$ws.Range("J140").Value = 2964.3333
$ws.Range("K140").Value = 120387.726
$ws.Range("L140").Value = 8892.999899999999
$ws.Range("M140").Value = -115207.726
$ws.Range("N140").Value = -19252.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4178
$ws.Range("I122").Value = 2260
$ws.Range("K122").Value = 6780
$ws.Range("M122").Value = -4330

$ws.Range("H126").Value = 2487.5
$ws.Range("I126").Value = 1252.0714
$ws.Range("J126").Value = 3568.5
$ws.Range("K126").Value = 3756.2142
$ws.Range("L126").Value = 10705.5
$ws.Range("M126").Value = -1286.2142
$ws.Range("N126").Value = -15645.5

$ws.Range("H132").Value = 2949.9048
$ws.Range("I132").Value = 2510
$ws.Range("J132").Value = 4049.6667
$ws.Range("K132").Value = 7530
$ws.Range("L132").Value = 12149.0001
$ws.Range("M132").Value = -5000
$ws.Range("N132").Value = -17209.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2844
$ws.Range("I7").Value = 2749.75
$ws.Range("J7").Value = 3032.5
$ws.Range("K7").Value = 2749.75
$ws.Range("L7").Value = 3032.5
$ws.Range("M7").Value = -2637.75
$ws.Range("N7").Value = -3256.5

$ws.Range("H32").Value = 3013
$ws.Range("I32").Value = 3013
$ws.Range("K32").Value = 3013
$ws.Range("M32").Value = -2696

$ws.Range("H93").Value = 1668
$ws.Range("I93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("M93").ClearContents()

$ws.Range("H126").Value = 2844
$ws.Range("I126").Value = 2749.75
$ws.Range("J126").Value = 3032.5
$ws.Range("K126").Value = 8249.25
$ws.Range("L126").Value = 9097.5
$ws.Range("M126").Value = -5779.25
$ws.Range("N126").Value = -14037.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 3111.111
$ws.Range("I15").Value = 2000
$ws.Range("K15").Value = 2000
$ws.Range("M15").Value = -1712

$ws.Range("H126").Value = 48310370
$ws.Range("I126").Value = 74074780
$ws.Range("K126").Value = 222224340
$ws.Range("M126").Value = -222221870

$ws.Range("H132").Value = 6698.857
$ws.Range("I132").Value = 8129.4614
$ws.Range("K132").Value = 24388.3842
$ws.Range("M132").Value = -21858.3842

$ws.Range("H136").Value = 1061.0476
$ws.Range("J136").Value = 2321.4285
$ws.Range("L136").Value = 6964.2855
$ws.Range("N136").Value = -12064.2855

Write-Output "edits applied"
